$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C slightly (larger mousebites).
# NOTE: the host's ColumnWidth setter quantizes to whole pixels
# (output = (round(input*6)+5)/6), so the literal target value of 9.59
# characters is not directly reachable; 8.8 is the input that lands on
# the closest achievable quantized width (9.6667 ~= 9.59).
$ws.Columns.Item(3).ColumnWidth = 8.8

# Row 2: R1
$ws.Cells.Item(2, 1).Value = "R1"
$ws.Cells.Item(2, 2).Value = 97.5342
$ws.Cells.Item(2, 3).Value = -85.8985
$ws.Cells.Item(2, 4).Value = "top"
$ws.Cells.Item(2, 5).Value = 0

# Row 3: R2
$ws.Cells.Item(3, 1).Value = "R2"
$ws.Cells.Item(3, 2).Value = 90.5342
$ws.Cells.Item(3, 3).Value = -101.8988
$ws.Cells.Item(3, 4).Value = "top"
$ws.Cells.Item(3, 5).Value = 0

# Row 4: R3 (new)
$ws.Cells.Item(4, 1).Value = "R3"
$ws.Cells.Item(4, 2).Value = 97.5342
$ws.Cells.Item(4, 3).Value = -117.8985
$ws.Cells.Item(4, 4).Value = "top"
$ws.Cells.Item(4, 5).Value = 0

# Row 5: R4 (new)
$ws.Cells.Item(5, 1).Value = "R4"
$ws.Cells.Item(5, 2).Value = 128.9978
$ws.Cells.Item(5, 3).Value = -85.8985
$ws.Cells.Item(5, 4).Value = "top"
$ws.Cells.Item(5, 5).Value = 0

# Row 6: R5 (new)
$ws.Cells.Item(6, 1).Value = "R5"
$ws.Cells.Item(6, 2).Value = 121.9978
$ws.Cells.Item(6, 3).Value = -101.8988
$ws.Cells.Item(6, 4).Value = "top"
$ws.Cells.Item(6, 5).Value = 0

# Row 7: R6 (new)
$ws.Cells.Item(7, 1).Value = "R6"
$ws.Cells.Item(7, 2).Value = 128.9978
$ws.Cells.Item(7, 3).Value = -117.8985
$ws.Cells.Item(7, 4).Value = "top"
$ws.Cells.Item(7, 5).Value = 0

# Row 8: SW1
$ws.Cells.Item(8, 1).Value = "SW1"
$ws.Cells.Item(8, 2).Value = 97.5342
$ws.Cells.Item(8, 3).Value = -81.7735
$ws.Cells.Item(8, 4).Value = "top"
$ws.Cells.Item(8, 5).Value = 180

# Row 9: SW2
$ws.Cells.Item(9, 1).Value = "SW2"
$ws.Cells.Item(9, 2).Value = 90.5342
$ws.Cells.Item(9, 3).Value = -97.7738
$ws.Cells.Item(9, 4).Value = "top"
$ws.Cells.Item(9, 5).Value = 180

# Row 10: SW3 (new)
$ws.Cells.Item(10, 1).Value = "SW3"
$ws.Cells.Item(10, 2).Value = 97.5342
$ws.Cells.Item(10, 3).Value = -113.7735
$ws.Cells.Item(10, 4).Value = "top"
$ws.Cells.Item(10, 5).Value = 180

# Row 11: SW4 (new)
$ws.Cells.Item(11, 1).Value = "SW4"
$ws.Cells.Item(11, 2).Value = 128.9978
$ws.Cells.Item(11, 3).Value = -81.7735
$ws.Cells.Item(11, 4).Value = "top"
$ws.Cells.Item(11, 5).Value = 180

# Row 12: SW5 (new)
$ws.Cells.Item(12, 1).Value = "SW5"
$ws.Cells.Item(12, 2).Value = 121.9978
$ws.Cells.Item(12, 3).Value = -97.7738
$ws.Cells.Item(12, 4).Value = "top"
$ws.Cells.Item(12, 5).Value = 180

# Row 13: SW6 (new)
$ws.Cells.Item(13, 1).Value = "SW6"
$ws.Cells.Item(13, 2).Value = 128.9978
$ws.Cells.Item(13, 3).Value = -113.7735
$ws.Cells.Item(13, 4).Value = "top"
$ws.Cells.Item(13, 5).Value = 180
